# GitHub Actions cryptos-list refresh: update Price (D) and Volume(1h) (E)
# for the coin rows (2-51) per the latest coinranking.com snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "41.525.71"
$ws.Range("E2").Value = "  +0.08%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.466.51"
$ws.Range("E3").Value = "  -0.43%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.61%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.01"
$ws.Range("E5").Value = "  -0.05%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "91.35"
$ws.Range("E6").Value = "  -1.75%  "

$ws.Range("E8").Value = "  -0.62%  "

$ws.Range("E9").Value = "  +3.84%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "32.53"
$ws.Range("E10").Value = "  -1.75%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0794"
$ws.Range("E11").Value = "  +1.84%  "

$ws.Range("E12").Value = "  +0.72%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.845.40"
$ws.Range("E13").Value = "  -0.44%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.85"
$ws.Range("E14").Value = "  -0.33%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.80"
$ws.Range("E15").Value = "  +2.95%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.457.40"
$ws.Range("E16").Value = "  -0.68%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.775"
$ws.Range("E17").Value = "  -1.19%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "41.534.33"
$ws.Range("E18").Value = "  +0.49%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.50"
$ws.Range("E19").Value = "  +3.33%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "70.98"
$ws.Range("E21").Value = "  +1.16%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.14"
$ws.Range("E22").Value = "  +0.77%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "238.18"
$ws.Range("E23").Value = "  +1.28%  "

$ws.Range("E24").Value = "  -0.80%  "

$ws.Range("E25").Value = "  -0.06%  "

$ws.Range("E26").Value = "  +1.18%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.55"
$ws.Range("E27").Value = "  +2.14%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.25"
$ws.Range("E28").Value = "  -0.07%  "

$ws.Range("E29").Value = "  -0.80%  "

$ws.Range("E30").Value = "  -3.21%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "156.23"
$ws.Range("E31").Value = "  +2.17%  "

$ws.Range("E32").Value = "  -0.17%  "

$ws.Range("E33").Value = "  +1.10%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0757"
$ws.Range("E34").Value = "  +0.69%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "17.23"
$ws.Range("E35").Value = "  -2.45%  "

$ws.Range("E36").Value = "  -7.74%  "

$ws.Range("E37").Value = "  -4.77%  "

$ws.Range("E38").Value = "  +1.09%  "

$ws.Range("E39").Value = "  +2.73%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.79"
$ws.Range("E40").Value = "  -3.50%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.98"
$ws.Range("E41").Value = "  -1.36%  "

$ws.Range("E42").Value = "  -0.87%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.941.89"
$ws.Range("E43").Value = "  -1.16%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0283"
$ws.Range("E44").Value = "  +0.29%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "18.79"
$ws.Range("E45").Value = "  -3.64%  "

$ws.Range("E46").Value = "  -2.16%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.03"
$ws.Range("E47").Value = "  +2.87%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.704.22"
$ws.Range("E48").Value = "  -0.63%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "96.67"
$ws.Range("E49").Value = "  +0.74%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "66.90"
$ws.Range("E50").Value = "  -2.21%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "52.14"
$ws.Range("E51").Value = "  +3.78%  "
